# Update the division-problem values in the practice worksheet table.
# Assigning to Table.Cell(r, c).Range.Text replaces the cell's visible
# text while Word preserves the end-of-cell mark, and addressing cells
# by (row, column) avoids any ambiguity from duplicate values (e.g. the
# two "77÷8=" cells in row 1) that a document-wide Find/Replace would hit.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1,1).Range.Text = "12÷3="
$t.Cell(1,2).Range.Text = "49÷4="
$t.Cell(1,3).Range.Text = "64÷4="
$t.Cell(1,4).Range.Text = "10÷3="
$t.Cell(1,5).Range.Text = "23÷4="

# Row 5
$t.Cell(5,1).Range.Text = "15÷8="
$t.Cell(5,2).Range.Text = "32÷4="
$t.Cell(5,3).Range.Text = "57÷9="
$t.Cell(5,4).Range.Text = "96÷8="
$t.Cell(5,5).Range.Text = "60÷8="

# Row 9
$t.Cell(9,1).Range.Text = "39÷5="
$t.Cell(9,2).Range.Text = "30÷7="
$t.Cell(9,3).Range.Text = "83÷9="
$t.Cell(9,4).Range.Text = "35÷4="
$t.Cell(9,5).Range.Text = "44÷3="

# Row 13 (column 2, "55÷2=", is unchanged by this edit)
$t.Cell(13,1).Range.Text = "23÷4="
$t.Cell(13,3).Range.Text = "73÷8="
$t.Cell(13,4).Range.Text = "92÷9="
$t.Cell(13,5).Range.Text = "64÷4="

# Row 17
$t.Cell(17,1).Range.Text = "15÷4="
$t.Cell(17,2).Range.Text = "71÷9="
$t.Cell(17,3).Range.Text = "65÷5="
$t.Cell(17,4).Range.Text = "63÷2="
$t.Cell(17,5).Range.Text = "57÷8="
